$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("J2").Value = 1.03
$ws.Range("K2").Value = 17
$ws.Range("L2").Value = 1.14
$ws.Range("M2").Value = 5.5
$ws.Range("N2").Value = 1.5
$ws.Range("O2").Value = 2.63
$ws.Range("J4").Value = 1.11
$ws.Range("K4").Value = 6.5
$ws.Range("J9").Value = 1.04
$ws.Range("K9").Value = 13
$ws.Range("N9").Value = 1.88
$ws.Range("O9").Value = 1.93
$ws.Range("G14").Value = 2.4
$ws.Range("I14").Value = 2.9
$ws.Range("V14").Value = 9.5
$ws.Range("G19").Value = 4.9
$ws.Range("H19").Value = 4.2
$ws.Range("I19").Value = 1.5
$ws.Range("N19").Value = 1.55
$ws.Range("O19").Value = 2.15
$ws.Range("T19").Value = 14
$ws.Range("U19").Value = 26
$ws.Range("Y19").Value = 35
$ws.Range("Z19").Value = 14.5
$ws.Range("AA19").Value = 7.4
$ws.Range("AB19").Value = 13
$ws.Range("AC19").Value = 45
$ws.Range("AD19").Value = 250
$ws.Range("AE19").Value = 7.3
$ws.Range("AF19").Value = 6.9
$ws.Range("AG19").Value = 7
$ws.Range("AH19").Value = 9.25
$ws.Range("AI19").Value = 9.5
$ws.Range("AJ19").Value = 17.5
$ws.Range("G31").Value = 2.38
$ws.Range("I31").Value = 3.1
$ws.Range("R31").Value = 1.63
$ws.Range("W31").Value = 23
$ws.Range("Z31").Value = 11
$ws.Range("AG31").Value = 11
$ws.Range("AH31").Value = 29
$ws.Range("AJ31").Value = 29
$ws.Range("S32").Value = 1.54
$ws.Range("R33").Value = 1.72
$ws.Range("R34").Value = 1.63
$ws.Range("R35").Value = 1.63
$ws.Range("R36").Value = 1.77
$ws.Range("S36").Value = 1.87
$ws.Range("G43").Value = 1.62
$ws.Range("I43").Value = 5
$ws.Range("AA43").Value = 8
$ws.Range("AB43").Value = 15
$ws.Range("G48").Value = 1.36
$ws.Range("I48").Value = 6.7
$ws.Range("N48").Value = 1.4
$ws.Range("O48").Value = 2.75
$ws.Range("U48").Value = 8.5
$ws.Range("X48").Value = 10.25
$ws.Range("Y48").Value = 19.5
$ws.Range("AE48").Value = 26
$ws.Range("AF48").Value = 50
$ws.Range("AG48").Value = 21
$ws.Range("AI48").Value = 60
$ws.Range("AJ48").Value = 45
